# Loading hosts from excel
# Adds a new "hosts" worksheet right after "vlan_definitions" and makes
# "bgp_neighbors" the active/selected sheet (matches the canonical diff).

$wb = $excel.ActiveWorkbook

# --- Insert the new "hosts" worksheet at position 2 -------------------------
$vlan = $wb.Worksheets.Item("vlan_definitions")
$physical = $wb.Worksheets.Item("physical_links")

$hosts = $wb.Worksheets.Add($physical)
$hosts.Name = "hosts"

# --- Populate the header row -------------------------------------------------
$hosts.Range("A1").Value = "use"
$hosts.Range("B1").Value = "name"
$hosts.Range("C1").Value = "platform"
$hosts.Range("D1").Value = "mgmt_host"
$hosts.Range("E1").Value = "tags"

# --- Populate the data row ----------------------------------------------------
$hosts.Range("A2").Value = $true
$hosts.Range("B2").Value = "TestHost-A"
$hosts.Range("C2").Value = "ios_xe"
$hosts.Range("D2").Value = "TestHost-A"
$hosts.Range("E2").Value = "one,two"

# --- Match column widths to content (bestFit) --------------------------------
$hosts.Columns.Item(1).ColumnWidth = 4.666666666666667
$hosts.Columns.Item(2).ColumnWidth = 9.833333333333334
$hosts.Columns.Item(3).ColumnWidth = 7.833333333333333
$hosts.Columns.Item(4).ColumnWidth = 10.166666666666666
$hosts.Columns.Item(5).ColumnWidth = 7.666666666666667

# --- Selection on the new sheet ------------------------------------------------
$hosts.Range("D8").Select() | Out-Null

# --- Move the active tab to "bgp_neighbors" -----------------------------------
$neighbors = $wb.Worksheets.Item("bgp_neighbors")
$neighbors.Activate() | Out-Null
$neighbors.Range("C2").Select() | Out-Null
